$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4490761131189629
$ws.Range("C2").Value = 0.04004426199074373
$ws.Range("D2").Value = 0.2875994519771581
$ws.Range("F2").Value = 1.498101536507988
$ws.Range("G2").Value = 0.002461461266048388
$ws.Range("J2").Value = 0.3501554756335281
$ws.Range("K2").Value = 0.4123360384060959
$ws.Range("N2").Value = 1.614479840121912
$ws.Range("O2").Value = 3.423728209994778
# Row 3
$ws.Range("B3").Value = 0.4108558859743141
$ws.Range("C3").Value = 0.03495564479990776
$ws.Range("D3").Value = 0.2773086194918193
$ws.Range("F3").Value = 1.497412064146481
$ws.Range("G3").Value = 0.00246396260673315
$ws.Range("J3").Value = 0.3388457038141723
$ws.Range("K3").Value = 0.3716892487181838
$ws.Range("N3").Value = 1.631882314868638
$ws.Range("O3").Value = 3.438086904713998
# Row 4
$ws.Range("B4").Value = 0.3875055162888543
$ws.Range("C4").Value = 0.03181911584162833
$ws.Range("D4").Value = 0.2711221100779113
$ws.Range("F4").Value = 1.497785075238632
$ws.Range("G4").Value = 0.002465580801574744
$ws.Range("J4").Value = 0.3321060681798542
$ws.Range("K4").Value = 0.3467961042615002
$ws.Range("N4").Value = 1.643111982848032
$ws.Range("O4").Value = 3.448702506425576
# Row 5
$ws.Range("B5").Value = 0.3780199714763626
$ws.Range("C5").Value = 0.03053797169705774
$ws.Range("D5").Value = 0.2686343921312755
$ws.Range("F5").Value = 1.498137376502292
$ws.Range("G5").Value = 0.002466261001336774
$ws.Range("J5").Value = 0.3294110401265726
$ws.Range("K5").Value = 0.3366685838594492
$ws.Range("N5").Value = 1.647825064109269
$ws.Range("O5").Value = 3.453480858761537
# Row 6
$ws.Range("B6").Value = 0.3764467258060336
$ws.Range("C6").Value = 0.03032506025050452
$ws.Range("D6").Value = 0.2682233252361925
$ws.Range("F6").Value = 1.498207974689727
$ws.Range("G6").Value = 0.002466375204218713
$ws.Range("J6").Value = 0.3289666382615479
$ws.Range("K6").Value = 0.3349879366684547
$ws.Range("N6").Value = 1.648615935772308
$ws.Range("O6").Value = 3.454301623741387
# Row 7
$ws.Range("B7").Value = 0.3873774691062124
$ws.Range("C7").Value = 0.03180184988332257
$ws.Range("D7").Value = 0.2710884247170497
$ws.Range("F7").Value = 1.497789015443544
$ws.Range("G7").Value = 0.00246558989090002
$ws.Range("J7").Value = 0.3320695138873333
$ws.Range("K7").Value = 0.3466594529399742
$ws.Range("N7").Value = 1.643174990823832
$ws.Range("O7").Value = 3.448765117231332
# Row 8
$ws.Range("B8").Value = 0.4358737814306721
$ws.Range("C8").Value = 0.03829225139955383
$ws.Range("D8").Value = 0.2840238209979589
$ws.Range("F8").Value = 1.497698548069522
$ws.Range("G8").Value = 0.00246230667139308
$ws.Range("J8").Value = 0.346213381598119
$ws.Range("K8").Value = 0.3983080027858819
$ws.Range("N8").Value = 1.620367200630064
$ws.Range("O8").Value = 3.428305611767144
# Row 9
$ws.Range("B9").Value = 0.5318862884136877
$ws.Range("C9").Value = 0.05092192898327141
$ws.Range("D9").Value = 0.3104347204198916
$ws.Range("F9").Value = 1.503839791265776
$ws.Range("G9").Value = 0.002456518972684715
$ws.Range("J9").Value = 0.3755758483228249
$ws.Range("K9").Value = 0.5000824528152918
$ws.Range("N9").Value = 1.579961292561864
$ws.Range("O9").Value = 3.402465838102671
# Row 10
$ws.Range("B10").Value = 0.6029667068298181
$ws.Range("C10").Value = 0.06013942043091447
$ws.Range("D10").Value = 0.330473170175253
$ws.Range("F10").Value = 1.512206796440125
$ws.Range("G10").Value = 0.002452659474522717
$ws.Range("J10").Value = 0.3981466537149032
$ws.Range("K10").Value = 0.5751403753493207
$ws.Range("N10").Value = 1.552907618536221
$ws.Range("O10").Value = 3.392197840763913
# Row 11
$ws.Range("B11").Value = 0.6354174256632916
$ws.Range("C11").Value = 0.06431900631814358
$ws.Range("D11").Value = 0.3397265393020916
$ws.Range("F11").Value = 1.516851108922978
$ws.Range("G11").Value = 0.002450988124911287
$ws.Range("J11").Value = 0.4086330448536017
$ws.Range("K11").Value = 0.6093450362159842
$ws.Range("N11").Value = 1.541171612212112
$ws.Range("O11").Value = 3.389421885269314
# Row 12
$ws.Range("B12").Value = 0.6477219239061469
$ws.Range("C12").Value = 0.06589971793742677
$ws.Range("D12").Value = 0.3432502763394041
$ws.Range("F12").Value = 1.518730324566476
$ws.Range("G12").Value = 0.002450367294980746
$ws.Range("J12").Value = 0.4126355026395458
$ws.Range("K12").Value = 0.6223057285193647
$ws.Range("N12").Value = 1.536809602930466
$ws.Range("O12").Value = 3.388643353157846
# Row 13
$ws.Range("B13").Value = 0.6450712220500066
$ws.Range("C13").Value = 0.06555937371265941
$ws.Range("D13").Value = 0.3424905029372098
$ws.Range("F13").Value = 1.518320241988533
$ws.Range("G13").Value = 0.002450500465784012
$ws.Range("J13").Value = 0.4117721008973803
$ws.Range("K13").Value = 0.6195140589017285
$ws.Range("N13").Value = 1.537745384294297
$ws.Range("O13").Value = 3.388798894671112
# Row 14
$ws.Range("B14").Value = 0.6364294043053178
$ws.Range("C14").Value = 0.06444909301336565
$ws.Range("D14").Value = 0.3400160454327477
$ws.Range("F14").Value = 1.517003297849456
$ws.Range("G14").Value = 0.002450936807043265
$ws.Range("J14").Value = 0.4089616982470545
$ws.Range("K14").Value = 0.6104111604566924
$ws.Range("N14").Value = 1.540811100003815
$ws.Range("O14").Value = 3.389352369480264
# Row 15
$ws.Range("B15").Value = 0.6311381238263891
$ws.Range("C15").Value = 0.0637687507189213
$ws.Range("D15").Value = 0.3385029288606063
$ws.Range("F15").Value = 1.516212326029162
$ws.Range("G15").Value = 0.00245120564963132
$ws.Range("J15").Value = 0.4072443469536182
$ws.Range("K15").Value = 0.604836414380344
$ws.Range("N15").Value = 1.542699641531293
$ws.Range("O15").Value = 3.389726902246991
# Row 16
$ws.Range("B16").Value = 0.6008482609871351
$ws.Range("C16").Value = 0.05986599742395526
$ws.Range("D16").Value = 0.3298712024802626
$ws.Range("F16").Value = 1.511920144357944
$ws.Range("G16").Value = 0.002452770393808031
$ws.Range("J16").Value = 0.3974657485820927
$ws.Range("K16").Value = 0.5729061884801752
$ws.Range("N16").Value = 1.553686087340454
$ws.Range("O16").Value = 3.392417402231047
# Row 17
$ws.Range("B17").Value = 0.5822957009148979
$ws.Range("C17").Value = 0.05746828075831445
$ws.Range("D17").Value = 0.3246111171098391
$ws.Range("F17").Value = 1.509501688766861
$ws.Range("G17").Value = 0.002453751878866369
$ws.Range("J17").Value = 0.391522956646881
$ws.Range("K17").Value = 0.5533331089110334
$ws.Range("N17").Value = 1.56057217545818
$ws.Range("O17").Value = 3.394553408954749
# Row 18
$ws.Range("B18").Value = 0.5716356989996996
$ws.Range("C18").Value = 0.056087911066939
$ws.Range("D18").Value = 0.3215986309503762
$ws.Range("F18").Value = 1.508189545197538
$ws.Range("G18").Value = 0.002454324346575199
$ws.Range("J18").Value = 0.3881254148126914
$ws.Range("K18").Value = 0.5420809132278634
$ws.Range("N18").Value = 1.56458659724328
$ws.Range("O18").Value = 3.395960340879498
# Row 19
$ws.Range("B19").Value = 0.5680283011603251
$ws.Range("C19").Value = 0.05562032609940104
$ws.Range("D19").Value = 0.3205808873188687
$ws.Range("F19").Value = 1.507758825206395
$ws.Range("G19").Value = 0.002454519540351412
$ws.Range("J19").Value = 0.386978601999715
$ws.Range("K19").Value = 0.5382721129621757
$ws.Range("N19").Value = 1.565955036757172
$ws.Range("O19").Value = 3.396467331984866
# Row 20
$ws.Range("B20").Value = 0.5842695243332798
$ws.Range("C20").Value = 0.05772365334665608
$ws.Range("D20").Value = 0.3251697201357615
$ws.Range("F20").Value = 1.50975097290484
$ws.Range("G20").Value = 0.002453646576407905
$ws.Range("J20").Value = 0.3921534453290718
$ws.Range("K20").Value = 0.555416108503465
$ws.Range("N20").Value = 1.559833578645581
$ws.Range("O20").Value = 3.394307566897226
# Row 21
$ws.Range("B21").Value = 0.6389672826622359
$ws.Range("C21").Value = 0.06477526435786274
$ws.Range("D21").Value = 0.3407423202852158
$ws.Range("F21").Value = 1.517386845558249
$ws.Range("G21").Value = 0.002450808315843955
$ws.Range("J21").Value = 0.4097863268458894
$ws.Range("K21").Value = 0.6130846862234876
$ws.Range("N21").Value = 1.539908394533268
$ws.Range("O21").Value = 3.389182399247943
# Row 22
$ws.Range("B22").Value = 0.6748091075093612
$ws.Range("C22").Value = 0.06937217181007327
$ws.Range("D22").Value = 0.3510346069913055
$ws.Range("F22").Value = 1.523079782658641
$ws.Range("G22").Value = 0.002449023699071744
$ws.Range("J22").Value = 0.42149398672872
$ws.Range("K22").Value = 0.6508215605513783
$ws.Range("N22").Value = 1.527365072405015
$ws.Range("O22").Value = 3.387422164472611
# Row 23
$ws.Range("B23").Value = 0.6556712606525252
$ws.Range("C23").Value = 0.0669198106687503
$ws.Range("D23").Value = 0.3455309670376039
$ws.Range("F23").Value = 1.519977082270117
$ws.Range("G23").Value = 0.002449969763994753
$ws.Range("J23").Value = 0.4152285857738747
$ws.Range("K23").Value = 0.630676553899832
$ws.Range("N23").Value = 1.534015831219417
$ws.Range("O23").Value = 3.388216154708743
# Row 24
$ws.Range("B24").Value = 0.5833771396863199
$ws.Range("C24").Value = 0.05760820528408317
$ws.Range("D24").Value = 0.3249171395452208
$ws.Range("F24").Value = 1.509638027772354
$ws.Range("G24").Value = 0.002453694158176443
$ws.Range("J24").Value = 0.3918683420669282
$ws.Range("K24").Value = 0.5544743823880651
$ws.Range("N24").Value = 1.560167325330822
$ws.Range("O24").Value = 3.394418154792049
# Row 25
$ws.Range("B25").Value = 0.5058164016550393
$ws.Range("C25").Value = 0.04751594264949688
$ws.Range("D25").Value = 0.3031782918682495
$ws.Range("F25").Value = 1.501501656544733
$ws.Range("G25").Value = 0.002458015448485164
$ws.Range("J25").Value = 0.3674577108974972
$ws.Range("K25").Value = 0.4724987284322424
$ws.Range("N25").Value = 1.590429790918178
$ws.Range("O25").Value = 3.407926221071648
